$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-typed cells (column D) to keep their original string type instead of
# being auto-coerced to numbers by Excel when the text looks numeric.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.926.62"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.14%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.919.80"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.90%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.72"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4565"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.64%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3799"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.43%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07737"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.36%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9757"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22.28"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.911.57"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.688"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.36%  "
$ws.Range("E14").Value = "  -0.23%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06975"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.86%  "
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "84.39"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.62%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009456"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.62"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.003"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "28.945.45"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.334"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.05"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.50%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.155.84"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.056"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "157.73"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.03"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.33%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.599"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.34%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "117.79"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.28%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.838"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.65%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09294"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8644"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.33%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.090"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.239"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.95%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.012"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05682"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.20%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.148"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.004"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02040"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.22%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.071"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +12.32%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.457"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.34%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5481"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.54%  "
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "9.312"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.000002762"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +14.52%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.157"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5158"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.65%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06941"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.76%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "11.12"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.32%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "110.55"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.60%  "
$ws.Range("E51").Value = "  -0.84%  "
